# clean up on BUH and Sales team roles
# Adds a new "LeadSoc" customer row (row 14) to Sheet1, mirroring the
# existing rows: cName / cEmail (mailto: hyperlink) / cUrl (http: hyperlink).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values
$ws.Range("A14").Value = "LeadSoc"
$ws.Range("B14").Value = "leadsoc@gmail.com"
$ws.Range("C14").Value = "www.leadsoc.com"

# Hook up hyperlinks for the email + url cells, same as every other row.
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:leadsoc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C14"), "http://www.leadsoc.com/")

# Hyperlinks.Add() stamps a freshly-minted "Hyperlink" style onto the
# cells instead of reusing the workbook's existing one; re-apply the
# formatting from the row above (which already carries the correct,
# pre-existing Hyperlink cell style) so B14:C14 line up with the rest
# of the table instead of picking up a duplicate style definition.
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
